# day3 Final Assignment - API automation run report update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Plain value edits on existing rows (Register / Login / Pagination blocks)
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Shaikh Test"
$ws.Range("B3").Value = "Shaikh Test"
$ws.Range("H6").Value = "1"

# ---------------------------------------------------------------------------
# 2) New "Adding_20_tasks" run row
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "Adding_20_tasks"
$ws.Range("B7").Value = "Shaikhh"
$ws.Range("D7").Value = "123456789"
$ws.Range("E7").Value = "34"

# ---------------------------------------------------------------------------
# 3) Insert the missing task17/17 pair into the task-name / task-note tables,
#    shifting the former R:T columns one slot right into S:U
#    (NOTE: ".Value" getter is unreliable on this host -- it can echo back
#    the property descriptor instead of the cell's contents -- so reads use
#    ".Value2" here; only the left-hand-side assignment uses ".Value".)
# ---------------------------------------------------------------------------
$ws.Range("U11").Value = $ws.Range("T11").Value2
$ws.Range("T11").Value = $ws.Range("S11").Value2
$ws.Range("S11").Value = $ws.Range("R11").Value2
$ws.Range("R11").Value = "task17"

$ws.Range("U12").Value = $ws.Range("T12").Value2
$ws.Range("T12").Value = $ws.Range("S12").Value2
$ws.Range("S12").Value = $ws.Range("R12").Value2
$ws.Range("R12").Value = "17"

# ---------------------------------------------------------------------------
# 4) New "Register_f" / "Login_f" Test_Func rows
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Register_f"
$ws.Range("B14").Value = "Test_Func"
$ws.Range("D14").Value = "123456789"
$ws.Range("E14").Value = "19"

$ws.Range("A15").Value = "Login_f"
$ws.Range("B15").Value = "Test_Func"
$ws.Range("D15").Value = "123456789"
$ws.Range("E15").Value = "33"

# ---------------------------------------------------------------------------
# 5) Hyperlinks: drop the old set (this clears every hyperlink on the sheet
#    in this host) and recreate all five, in left-to-right / top-to-bottom
#    order, so relationship ids come out rId1..rId5 the way Excel would
#    assign them.
# ---------------------------------------------------------------------------
$ws.Range("C2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:testhash39@test.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:testhash39@test.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:taprt29816@test.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C14"), "mailto:register@test.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C15"), "mailto:register@test.com") | Out-Null

$ws.Range("C2").Value = "testhash39@test.com"
$ws.Range("C3").Value = "testhash39@test.com"
$ws.Range("C7").Value = "taprt29816@test.com"
$ws.Range("C14").Value = "register@test.com"
$ws.Range("C15").Value = "register@test.com"

# Restore the standard hyperlink-cell formatting (style "1" in the original
# workbook) -- Hyperlinks.Add() stamps its own new style, so copy the look
# back in from an untouched hyperlink-styled cell (B4 always stays empty).
$ws.Range("B4").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6) Selection cursor ends on B7, matching the saved workbook state
# ---------------------------------------------------------------------------
$ws.Range("B7").Select() | Out-Null

Write-Output "edit complete"
